$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-26 Saturday" "2025-07-27 Sunday"

Replace-Text "68÷2=34, 0" "46÷2=23, 0"
Replace-Text "85÷4=21, 1" "55÷2=27, 1"
Replace-Text "44÷3=14, 2" "29÷6=4, 5"
Replace-Text "29÷2=14, 1" "81÷4=20, 1"
Replace-Text "63÷5=12, 3" "76÷3=25, 1"

Replace-Text "48÷8=6, 0" "21÷5=4, 1"
Replace-Text "11÷7=1, 4" "12÷3=4, 0"
Replace-Text "50÷9=5, 5" "96÷8=12, 0"
Replace-Text "41÷2=20, 1" "60÷3=20, 0"
Replace-Text "13÷5=2, 3" "43÷6=7, 1"

Replace-Text "76÷8=9, 4" "42÷6=7, 0"
Replace-Text "29÷9=3, 2" "89÷4=22, 1"
Replace-Text "38÷9=4, 2" "74÷3=24, 2"
Replace-Text "22÷6=3, 4" "89÷5=17, 4"
Replace-Text "84÷8=10, 4" "11÷3=3, 2"

Replace-Text "64÷7=9, 1" "12÷7=1, 5"
Replace-Text "27÷7=3, 6" "66÷6=11, 0"
Replace-Text "87÷5=17, 2" "15÷9=1, 6"
Replace-Text "11÷9=1, 2" "63÷5=12, 3"
Replace-Text "67÷8=8, 3" "70÷4=17, 2"

Replace-Text "39÷8=4, 7" "44÷4=11, 0"
Replace-Text "10÷7=1, 3" "63÷5=12, 3"
Replace-Text "64÷8=8, 0" "47÷5=9, 2"
Replace-Text "36÷5=7, 1" "88÷9=9, 7"
Replace-Text "53÷6=8, 5" "48÷9=5, 3"
